# Insert two new price-record rows at the top of the data block that starts
# at row 326 (this pushes the former rows 326-362 down to rows 328-364).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(326).Insert()
$ws.Rows.Item(326).Insert()

# New row 326: Sutil De Gase, Primera
$ws.Cells.Item(326,1).Value  = 1
$ws.Cells.Item(326,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(326,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(326,4).Value  = 44946
$ws.Cells.Item(326,5).Value  = 15
$ws.Cells.Item(326,6).Value  = "Fruta"
$ws.Cells.Item(326,7).Value  = 100102
$ws.Cells.Item(326,8).Value  = "Cítricos"
$ws.Cells.Item(326,9).Value  = 100102003
$ws.Cells.Item(326,10).Value = "Limón"
$ws.Cells.Item(326,11).Value = "Sutil De Gase"
$ws.Cells.Item(326,12).Value = "Primera"
$ws.Cells.Item(326,13).Value = 380
$ws.Cells.Item(326,14).Value = 26000
$ws.Cells.Item(326,15).Value = 27000
$ws.Cells.Item(326,16).Value = 26658
$ws.Cells.Item(326,17).Value = "$/caja 24 kilos"
$ws.Cells.Item(326,18).Value = "Perú"
$ws.Cells.Item(326,19).Value = 1111
$ws.Cells.Item(326,20).Value = 24

# New row 327: Tahití, Primera
$ws.Cells.Item(327,1).Value  = 1
$ws.Cells.Item(327,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(327,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(327,4).Value  = 44946
$ws.Cells.Item(327,5).Value  = 15
$ws.Cells.Item(327,6).Value  = "Fruta"
$ws.Cells.Item(327,7).Value  = 100102
$ws.Cells.Item(327,8).Value  = "Cítricos"
$ws.Cells.Item(327,9).Value  = 100102003
$ws.Cells.Item(327,10).Value = "Limón"
$ws.Cells.Item(327,11).Value = "Tahití"
$ws.Cells.Item(327,12).Value = "Primera"
$ws.Cells.Item(327,13).Value = 700
$ws.Cells.Item(327,14).Value = 30000
$ws.Cells.Item(327,15).Value = 31000
$ws.Cells.Item(327,16).Value = 30571
$ws.Cells.Item(327,17).Value = "$/caja 24 kilos"
$ws.Cells.Item(327,18).Value = "Perú"
$ws.Cells.Item(327,19).Value = 1274
$ws.Cells.Item(327,20).Value = 24
